# Inventory update: the item "L8SEY6 / Almohadilla+Chip Epson T6712" (which
# was appended as the very last row of the sheet, row 86) actually belongs
# right after the "Almohadilla..." group, before "Buje de rodillo superior
# para Kyocera". This moves it up to row 12, pushing every row that was at
# 12-85 down by one (to 13-86).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 86
$targetRow = 12

# 1) Insert a new blank row at the target position; this shifts the old
#    rows $targetRow..$lastRow down by one (to $targetRow+1..$lastRow+1).
$ws.Rows.Item($targetRow).Insert()

# 2) The row that used to be last is now one row further down (it got
#    shifted too). Copy its values into the freshly-inserted blank row...
$movedRow = $lastRow + 1
$ws.Range("A" + $movedRow + ":J" + $movedRow).Copy()
$ws.Range("A" + $targetRow).PasteSpecial()

# 3) ...then delete the now-duplicated trailing row, restoring the sheet to
#    its original row count with the item relocated to row 12.
$ws.Rows.Item($movedRow).Delete()

# PasteSpecial copies values, not formulas, so re-establish the H/I
# formulas for the relocated row (same pattern as every other data row).
$ws.Range("H" + $targetRow).Formula = "=(E" + $targetRow + "-D" + $targetRow + ")*G" + $targetRow
$ws.Range("I" + $targetRow).Formula = "=D" + $targetRow + "*F" + $targetRow
